$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$paths = @(
    "D:\AProg\html\PortfolioC\competences.html",
    "D:\AProg\html\PortfolioC\experiences.html",
    "D:\AProg\html\PortfolioC\loisirs.html",
    "D:\AProg\html\PortfolioC\img\testfile\test\Nouveau document texte - Copie (2).txt",
    "D:\AProg\html\PortfolioC\img\testfile\test\Nouveau document texte - Copie (3).txt",
    "D:\AProg\html\PortfolioC\img\testfile\test\Nouveau document texte - Copie (4).txt",
    "D:\AProg\html\PortfolioC\img\testfile\test\Nouveau document texte - Copie (5).txt",
    "D:\AProg\html\PortfolioC\img\testfile\test\Nouveau document texte - Copie.txt",
    "D:\AProg\html\PortfolioC\img\testfile\test\Nouveau document texte.txt",
    "D:\AProg\html\PortfolioC\img\testfile\test\tezst\Nouveau document texte - Copie (2).txt",
    "D:\AProg\html\PortfolioC\img\testfile\test\tezst\Nouveau document texte - Copie (3).txt",
    "D:\AProg\html\PortfolioC\img\testfile\test\tezst\Nouveau document texte - Copie.txt",
    "D:\AProg\html\PortfolioC\img\testfile\test\tezst\Nouveau document texte.txt"
)

$row = 2
foreach ($p in $paths) {
    $ws.Cells.Item($row, 1).Value = $p
    $ws.Cells.Item($row, 2).Value = "None"
    $ws.Cells.Item($row, 3).Value = "None"
    $ws.Cells.Item($row, 4).Value = "File not found"
    $row++
}
